$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A130:A131").NumberFormat = "@"

$ws.Range("A130").Value = "2022/01/10"
$ws.Range("B130").Value = 3204.4
$ws.Range("C130").Value = 3209.1
$ws.Range("D130").Value = 1.24
$ws.Range("E130").Value = 1.24

$ws.Range("A131").Value = "2022/01/12"
$ws.Range("B131").Value = 3615.9
$ws.Range("C131").Value = 3615.3
$ws.Range("D131").Value = 1.23
$ws.Range("E131").Value = 1.23
